$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the original column A (the TAXON-numbering column with values 4/11/14/18
# and the bold centered style) -- everything else shifts one column left
# (B->A, C->B, D->C, E->D, F->E), matching the diff's dimension change
# from A1:F5 to A1:E5 and the column-shifted cell contents.
$ws.Columns("A").Delete()

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION" (now sitting
# in column D after the shift) by editing out the underscore character.
$ws.Range("D1").Characters(6, 1).Text = ""
